$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 1529
$ws.Range("I2").Value = 1145.5714
$ws.Range("J2").Value = 2200
$ws.Range("K2").Value = 1145.5714
$ws.Range("L2").Value = 2200
$ws.Range("M2").Value = -1032.5714
$ws.Range("N2").Value = -2426

# Row 5
$ws.Range("H5").Value = 40.714287
$ws.Range("I5").Value = 40.714287
$ws.Range("K5").Value = 40.714287
$ws.Range("M5").Value = 74.285713

# Row 9
$ws.Range("H9").Value = 88.25
$ws.Range("I9").Value = 65
$ws.Range("J9").Value = 158
$ws.Range("K9").Value = 65
$ws.Range("L9").Value = 158
$ws.Range("M9").Value = 104
$ws.Range("N9").Value = -496

# Row 33
$ws.Range("H33").Value = 1338
$ws.Range("I33").Value = 1172.5
$ws.Range("K33").Value = 1172.5
$ws.Range("M33").Value = -943.5

# Row 41
$ws.Range("H41").Value = 1587.091
$ws.Range("I41").Value = 1468.6666
$ws.Range("J41").Value = 1729.2
$ws.Range("K41").Value = 1468.6666
$ws.Range("L41").Value = 1729.2
$ws.Range("M41").Value = -1028.6666
$ws.Range("N41").Value = -2609.2

# Row 52
$ws.Range("H52").Value = 1000
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 1000
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 3000
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value = -3320

# Row 53
$ws.Range("H53").Value = 427.42105
$ws.Range("I53").Value = 498.6154
$ws.Range("J53").Value = 273.16666
$ws.Range("K53").Value = 498.6154
$ws.Range("L53").Value = 273.16666
$ws.Range("M53").Value = 138.3846
$ws.Range("N53").Value = -1547.16666

# Row 62
$ws.Range("H62").Value = 5643.125
$ws.Range("J62").Value = 7997.5
$ws.Range("L62").Value = 7997.5
$ws.Range("N62").Value = -9245.5

# Row 64
$ws.Range("H64").Value = 2250
$ws.Range("I64").Value = 2000
$ws.Range("K64").Value = 2000
$ws.Range("M64").Value = -1752

# Row 65
$ws.Range("H65").Value = 5643.125
$ws.Range("J65").Value = 7997.5
$ws.Range("L65").Value = 39987.5
$ws.Range("N65").Value = -46227.5

# Row 67
$ws.Range("H67").Value = 2250
$ws.Range("I67").Value = 2000
$ws.Range("K67").Value = 2000
$ws.Range("M67").Value = -1142

# Row 74
$ws.Range("H74").Value = 253375
$ws.Range("I74").Value = 6750
$ws.Range("K74").Value = 6750
$ws.Range("M74").Value = -5814

# Row 77
$ws.Range("H77").Value = 253375
$ws.Range("I77").Value = 6750
$ws.Range("K77").Value = 33750
$ws.Range("M77").Value = -29070

# Row 100
$ws.Range("H100").Value = 1976.8572
$ws.Range("I100").Value = 2627.8
$ws.Range("K100").Value = 2627.8
$ws.Range("M100").Value = -2086.8

# Row 107
$ws.Range("H107").Value = 1343.5454
$ws.Range("I107").Value = 697.6539
$ws.Range("J107").Value = 3742.5715
$ws.Range("K107").Value = 697.6539
$ws.Range("L107").Value = 3742.5715
$ws.Range("M107").Value = 1222.3461
$ws.Range("N107").Value = -7582.5715

# Row 137
$ws.Range("H137").Value = 2122.2856
$ws.Range("I137").Value = 1619.3636
$ws.Range("J137").Value = 3966.3333
$ws.Range("K137").Value = 4858.0908
$ws.Range("L137").Value = 11898.9999
$ws.Range("M137").Value = -2308.0908
$ws.Range("N137").Value = -16998.9999

# Row 138
$ws.Range("H138").Value = 3937.36
$ws.Range("J138").Value = 3951.625
$ws.Range("L138").Value = 11854.875
$ws.Range("N138").Value = -22134.875

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 12108.119
$ws.Range("I32").Value = 10054.923
$ws.Range("K32").Value = 10054.923
$ws.Range("M32").Value = -9767.923000000001

# Row 45
$ws.Range("H45").Value = 500
$ws.Range("I45").Value = 500
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 500
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -123
$ws.Range("N45").ClearContents()

# Row 132
$ws.Range("H132").Value = 937.3333
$ws.Range("I132").Value = 937.3333
$ws.Range("K132").Value = 2811.9999
$ws.Range("M132").Value = -281.9998999999998

$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 288.83334
$ws.Range("I22").Value = 266.5
$ws.Range("K22").Value = 266.5
$ws.Range("M22").Value = -93.5

# Row 94
$ws.Range("H94").Value = 1583.909
$ws.Range("I94").Value = 1341.3
$ws.Range("J94").Value = 4010
$ws.Range("K94").Value = 1341.3
$ws.Range("L94").Value = 4010
$ws.Range("M94").Value = -890.3
$ws.Range("N94").Value = -4912

# Row 134
$ws.Range("H134").Value = 11735.714
$ws.Range("I134").Value = 8652.777
$ws.Range("J134").Value = 15000
$ws.Range("K134").Value = 25958.331
$ws.Range("L134").Value = 45000
$ws.Range("M134").Value = -23423.331
$ws.Range("N134").Value = -50070

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 894
$ws.Range("J22").Value = 895.3333
$ws.Range("L22").Value = 895.3333
$ws.Range("N22").Value = -1595.3333

# Row 31
$ws.Range("H31").Value = 3542.375
$ws.Range("I31").Value = 2756.3
$ws.Range("J31").Value = 4852.5
$ws.Range("K31").Value = 2756.3
$ws.Range("L31").Value = 4852.5
$ws.Range("M31").Value = -2461.3
$ws.Range("N31").Value = -5442.5

# Row 34
$ws.Range("H34").Value = 3542.375
$ws.Range("I34").Value = 2756.3
$ws.Range("J34").Value = 4852.5
$ws.Range("K34").Value = 2756.3
$ws.Range("L34").Value = 4852.5
$ws.Range("M34").Value = -2554.3
$ws.Range("N34").Value = -5256.5

# Row 64
$ws.Range("H64").Value = 30000
$ws.Range("J64").Value = 30000
$ws.Range("L64").Value = 30000
$ws.Range("N64").Value = -30496

# Row 67
$ws.Range("H67").Value = 30000
$ws.Range("J67").Value = 30000
$ws.Range("L67").Value = 30000
$ws.Range("N67").Value = -31716

# Row 105
$ws.Range("H105").Value = 1933.6364
$ws.Range("I105").Value = 1877
$ws.Range("K105").Value = 1877
$ws.Range("M105").Value = -130

# Row 132
$ws.Range("H132").Value = 1302
$ws.Range("I132").Value = 1302
$ws.Range("K132").Value = 3906
$ws.Range("M132").Value = -1376

# Row 134
$ws.Range("H134").Value = 4370.7856
$ws.Range("I134").Value = 4432.9165
$ws.Range("K134").Value = 13298.7495
$ws.Range("M134").Value = -10763.7495

$ws = $wb.Worksheets.Item("CUL")
# Row 26
$ws.Range("H26").Value = 5053.4546
$ws.Range("J26").Value = 6043.1113
$ws.Range("L26").Value = 18129.3339
$ws.Range("N26").Value = -18705.3339

# Row 131
$ws.Range("H131").Value = 2346.5454
$ws.Range("J131").Value = 3000
$ws.Range("L131").Value = 9000
$ws.Range("N131").Value = -19080

# Row 139
$ws.Range("H139").Value = 1150.1666
$ws.Range("I139").Value = 1150.1666
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 3450.4998
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = 1689.5002
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 93
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("N93").ClearContents()

# Row 97
$ws.Range("H97").Value = 8000
$ws.Range("J97").Value = 8000
$ws.Range("L97").Value = 8000
$ws.Range("N97").Value = -9982

$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 1867.875
$ws.Range("I113").Value = 975.5
$ws.Range("K113").Value = 2926.5
$ws.Range("M113").Value = -756.5

# Row 130
$ws.Range("H130").Value = 49332.668
$ws.Range("J130").Value = 49332.668
$ws.Range("L130").Value = 49332.668
$ws.Range("N130").Value = -59372.668

# Row 132
$ws.Range("H132").Value = 3999
$ws.Range("I132").Value = 3999
$ws.Range("K132").Value = 11997
$ws.Range("M132").Value = -9467

# Row 136
$ws.Range("H136").Value = 3947.8333
$ws.Range("I136").Value = 3852.182
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 11556.546
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -9006.545999999998
$ws.Range("N136").Value = -20100
